$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarter period labels (shift to the next window)
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish-date labels (shift to the next window)
$ws.Range("D9").Value = "1400-10-29 (3)"
$ws.Range("E9").Value = "1401-04-26 (9)"
$ws.Range("F9").Value = "1401-04-30 (2)"
$ws.Range("G9").Value = "1401-09-06 (4)"
$ws.Range("H9").Value = "1401-11-01 (4)"
$ws.Range("I9").Value = "1402-02-30 (8)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-09-06 (2)"
$ws.Range("L9").Value = "1401-11-01 (2)"
$ws.Range("M9").Value = "1402-02-30"

# Data rows 12-56: shift quarterly figures left by one column, append newest quarter at M
$ws.Range("D12").Value = 705417
$ws.Range("E12").Value = -494694
$ws.Range("F12").Value = 457
$ws.Range("G12").Value = 1013989
$ws.Range("H12").Value = 80283
$ws.Range("I12").Value = 255753
$ws.Range("J12").Value = 226591
$ws.Range("K12").Value = 189185
$ws.Range("L12").Value = 263947
$ws.Range("M12").Value = 560913

$ws.Range("D13").Value = -31693
$ws.Range("E13").Value = -23780
$ws.Range("F13").Value = -23769
$ws.Range("G13").Value = -27438
$ws.Range("H13").Value = -36651
$ws.Range("I13").Value = -52189
$ws.Range("J13").Value = -28890
$ws.Range("K13").Value = -42338
$ws.Range("L13").Value = -45623
$ws.Range("M13").Value = -52328

$ws.Range("D14").Value = 673724
$ws.Range("E14").Value = -518474
$ws.Range("F14").Value = -23312
$ws.Range("G14").Value = 986551
$ws.Range("H14").Value = 43632
$ws.Range("I14").Value = 203564
$ws.Range("J14").Value = 197701
$ws.Range("K14").Value = 146847
$ws.Range("L14").Value = 218324
$ws.Range("M14").Value = 508585


$ws.Range("D16").Value = 725
$ws.Range("E16").Value = 2809
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0

$ws.Range("D17").Value = -9135
$ws.Range("E17").Value = -66656
$ws.Range("F17").Value = -16343
$ws.Range("G17").Value = -63965
$ws.Range("H17").Value = -31320
$ws.Range("I17").Value = -55781
$ws.Range("J17").Value = -44802
$ws.Range("K17").Value = -52206
$ws.Range("L17").Value = -51611
$ws.Range("M17").Value = -66728

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0

$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 71
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0

$ws.Range("D20").Value = 0
$ws.Range("E20").Value = -62
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = -500
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = -762
$ws.Range("L20").Value = -693
$ws.Range("M20").Value = -3494

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0

$ws.Range("D26").Value = -165000
$ws.Range("E26").Value = 165000
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -198453

$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0

$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0

$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0

$ws.Range("D31").Value = 1895
$ws.Range("E31").Value = 4727
$ws.Range("F31").Value = 586
$ws.Range("G31").Value = 2212
$ws.Range("H31").Value = -2798
$ws.Range("I31").Value = 8868
$ws.Range("J31").Value = 3407
$ws.Range("K31").Value = 14332
$ws.Range("L31").Value = 2831
$ws.Range("M31").Value = 1702

$ws.Range("D32").Value = -171515
$ws.Range("E32").Value = 105889
$ws.Range("F32").Value = -15757
$ws.Range("G32").Value = -61753
$ws.Range("H32").Value = -34118
$ws.Range("I32").Value = -40746
$ws.Range("J32").Value = -41395
$ws.Range("K32").Value = -38636
$ws.Range("L32").Value = -49473
$ws.Range("M32").Value = -266973

$ws.Range("D33").Value = 502209
$ws.Range("E33").Value = -412585
$ws.Range("F33").Value = -39069
$ws.Range("G33").Value = 924798
$ws.Range("H33").Value = 9514
$ws.Range("I33").Value = 162818
$ws.Range("J33").Value = 156306
$ws.Range("K33").Value = 108211
$ws.Range("L33").Value = 168851
$ws.Range("M33").Value = 241612


$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0

$ws.Range("D36").Value = 52
$ws.Range("E36").Value = 52
$ws.Range("F36").Value = 52
$ws.Range("G36").Value = 52
$ws.Range("H36").Value = 52
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 52
$ws.Range("K36").Value = 52
$ws.Range("L36").Value = 52
$ws.Range("M36").Value = 0

$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0

$ws.Range("D38").Value = 0
$ws.Range("E38").Value = -165000
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -21000

$ws.Range("D39").Value = 150000
$ws.Range("E39").Value = 295000
$ws.Range("F39").Value = 100000
$ws.Range("G39").Value = 193000
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 7000
$ws.Range("J39").Value = 280000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 425313

$ws.Range("D40").Value = -47234
$ws.Range("E40").Value = -300318
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = -268374
$ws.Range("H40").Value = -123374
$ws.Range("I40").Value = -53252
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = -65448
$ws.Range("L40").Value = -68495
$ws.Range("M40").Value = -71370

$ws.Range("D41").Value = -1689
$ws.Range("E41").Value = -8873
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = -22788
$ws.Range("H41").Value = -5053
$ws.Range("I41").Value = -3658
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = -12600
$ws.Range("L41").Value = -9655
$ws.Range("M41").Value = -6577

$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0

$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0

$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0

$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0

$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0

$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0

$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0

$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0

$ws.Range("D50").Value = -596484
$ws.Range("E50").Value = 214357
$ws.Range("F50").Value = -14320
$ws.Range("G50").Value = -214533
$ws.Range("H50").Value = -150432
$ws.Range("I50").Value = -283214
$ws.Range("J50").Value = -1129
$ws.Range("K50").Value = -289130
$ws.Range("L50").Value = -232469
$ws.Range("M50").Value = -493189

$ws.Range("D51").Value = -495407
$ws.Range("E51").Value = 35166
$ws.Range("F51").Value = 85680
$ws.Range("G51").Value = -312695
$ws.Range("H51").Value = -278859
$ws.Range("I51").Value = -333124
$ws.Range("J51").Value = 278871
$ws.Range("K51").Value = -367178
$ws.Range("L51").Value = -310619
$ws.Range("M51").Value = -166823

$ws.Range("D52").Value = 6802
$ws.Range("E52").Value = -377419
$ws.Range("F52").Value = 46611
$ws.Range("G52").Value = 612103
$ws.Range("H52").Value = -269345
$ws.Range("I52").Value = -170306
$ws.Range("J52").Value = 435177
$ws.Range("K52").Value = -258967
$ws.Range("L52").Value = -141768
$ws.Range("M52").Value = 74789

$ws.Range("D53").Value = 421995
$ws.Range("E53").Value = 427461
$ws.Range("F53").Value = 71397
$ws.Range("G53").Value = 118008
$ws.Range("H53").Value = 731828
$ws.Range("I53").Value = 460766
$ws.Range("J53").Value = 290385
$ws.Range("K53").Value = 725562
$ws.Range("L53").Value = 467552
$ws.Range("M53").Value = 323870

$ws.Range("D54").Value = -1336
$ws.Range("E54").Value = 21355
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1717
$ws.Range("H54").Value = -1717
$ws.Range("I54").Value = -75
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 957
$ws.Range("L54").Value = -1914
$ws.Range("M54").Value = -440

$ws.Range("D55").Value = 427461
$ws.Range("E55").Value = 71397
$ws.Range("F55").Value = 118008
$ws.Range("G55").Value = 731828
$ws.Range("H55").Value = 460766
$ws.Range("I55").Value = 290385
$ws.Range("J55").Value = 725562
$ws.Range("K55").Value = 467552
$ws.Range("L55").Value = 323870
$ws.Range("M55").Value = 398219

$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
